$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original formatting/style of the data range so re-typing
# numeric-looking text (e.g. "245.07") doesn't get auto-coerced into a
# floating point number and lose its exact textual representation.
$dataRange = $ws.Range("B2:E51")
$origStyle = $dataRange.Style
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '30.337.96'
$ws.Range('E2').Value = '  -0.01%  '
$ws.Range('D3').Value = '1.934.87'
$ws.Range('E3').Value = '  +0.18%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '0.7555'
$ws.Range('E5').Value = '  +5.83%  '
$ws.Range('D6').Value = '245.07'
$ws.Range('E6').Value = '  -2.46%  '
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  -2.44%  '
$ws.Range('D9').Value = '27.65'
$ws.Range('E9').Value = '  +1.31%  '
$ws.Range('E10').Value = '  -2.45%  '
$ws.Range('D11').Value = '0.7780'
$ws.Range('E11').Value = '  -2.33%  '
$ws.Range('D12').Value = '0.08009'
$ws.Range('E12').Value = '  -0.88%  '
$ws.Range('D13').Value = '1.932.42'
$ws.Range('E13').Value = '  +0.10%  '
$ws.Range('D14').Value = '5.345'
$ws.Range('E14').Value = '  -1.22%  '
$ws.Range('E15').Value = '  -0.31%  '
$ws.Range('E16').Value = '  -2.51%  '
$ws.Range('D17').Value = '30.349.39'
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('D18').Value = '252.55'
$ws.Range('E18').Value = '  +0.44%  '
$ws.Range('D19').Value = '0.000007922'
$ws.Range('E19').Value = '  -2.35%  '
$ws.Range('D20').Value = '5.756'
$ws.Range('E20').Value = '  -0.29%  '
$ws.Range('D21').Value = '2.183.24'
$ws.Range('E21').Value = '  +0.08%  '
$ws.Range('D22').Value = '0.9997'
$ws.Range('E22').Value = '  -0.12%  '
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').Value = '  -0.23%  '
$ws.Range('D24').Value = '6.667'
$ws.Range('E24').Value = '  -3.42%  '
$ws.Range('D25').Value = '9.460'
$ws.Range('E25').Value = '  -2.45%  '
$ws.Range('E26').Value = '  +0.41%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '18.97'
$ws.Range('E27').Value = '  -1.13%  '
$ws.Range('B28').Value = 'Stellar'
$ws.Range('C28').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D28').Value = '0.1328'
$ws.Range('E28').Value = '  +3.99%  '
$ws.Range('D29').Value = '2.176'
$ws.Range('E29').Value = '  -5.97%  '
$ws.Range('D30').Value = '1.366'
$ws.Range('E30').Value = '  -0.05%  '
$ws.Range('D31').Value = '1.517'
$ws.Range('E31').Value = '  -1.70%  '
$ws.Range('D32').Value = '4.388'
$ws.Range('E32').Value = '  -0.74%  '
$ws.Range('D33').Value = '4.121'
$ws.Range('E33').Value = '  -1.61%  '
$ws.Range('D34').Value = '0.05150'
$ws.Range('E34').Value = '  -0.88%  '
$ws.Range('E35').Value = '  +1.38%  '
$ws.Range('D36').Value = '0.7501'
$ws.Range('E36').Value = '  +0.61%  '
$ws.Range('D37').Value = '2.769'
$ws.Range('E37').Value = '  +0.36%  '
$ws.Range('D38').Value = '0.01960'
$ws.Range('E38').Value = '  +0.36%  '
$ws.Range('D39').Value = '2.800'
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('D40').Value = '77.48'
$ws.Range('E40').Value = '  -1.79%  '
$ws.Range('D41').Value = '6.399'
$ws.Range('E41').Value = '  -0.31%  '
$ws.Range('D42').Value = '0.4449'
$ws.Range('E42').Value = '  -1.28%  '
$ws.Range('D43').Value = '1.963'
$ws.Range('E43').Value = '  -2.94%  '
$ws.Range('E44').Value = '  -0.03%  '
$ws.Range('D45').Value = '0.8340'
$ws.Range('E45').Value = '  -0.66%  '
$ws.Range('D46').Value = '100.60'
$ws.Range('E46').Value = '  -1.05%  '
$ws.Range('D47').Value = '9.750'
$ws.Range('E47').Value = '  -0.22%  '
$ws.Range('D48').Value = '7.465'
$ws.Range('E48').Value = '  +1.04%  '
$ws.Range('D49').Value = '37.45'
$ws.Range('E49').Value = '  +2.73%  '
$ws.Range('D50').Value = '978.97'
$ws.Range('E50').Value = '  +11.22%  '
$ws.Range('D51').Value = '0.06012'
$ws.Range('E51').Value = '  -1.27%  '

# Row 27 / Row 28 content swap (Stellar <-> EthereumClassic), new values
# for B/C already applied above alongside D/E.

$dataRange.Style = $origStyle

Write-Host "cryptos sheet updated"
